$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 'Infection_and_sepsis/SARI/Management of patients with SARI-additional Information.pdf'
$ws.Range("A21").Value = 'Infection_and_sepsis/SARI/Suspected Influeza A-H7N9 Guideline.pdf'
$ws.Range("A23").Value = 'Infection_and_sepsis/SARI/Management of Patients with severe acute respiratory infection SARI.pdf'
$ws.Range("A29").Value = 'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf'
$ws.Range("A30").Value = 'End_of_life_care/Reasons to report a death to PF.pdf'
$ws.Range("A31").Value = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf'
$ws.Range("A32").Value = 'Neurological/SOP -  Femoral site care.pdf'
$ws.Range("A38").Value = 'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf'
$ws.Range("A39").Value = 'Covid-19/WGH/CoVid intubation checklist WGH.pdf'
$ws.Range("A40").Value = 'Airway/Emergency intubation checklist_em_pub.pdf'
$ws.Range("A41").Value = 'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf'
$ws.Range("A42").Value = 'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf'
$ws.Range("A50").Value = 'Delirium/Managing a Potentially Violent Patient.pdf'
$ws.Range("A51").Value = 'Delirium/Risk assessment posi mit.pdf'
$ws.Range("A52").Value = 'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf'
$ws.Range("A53").Value = 'Breathing(Respiratory)/HFNO.pdf'
$ws.Range("A54").Value = 'Drugs/ketamine_in_asthma.pdf'
$ws.Range("A55").Value = 'Delirium/Drugs Causing Delirium and Agitiation.pdf'
$ws.Range("A57").Value = 'Airway/McGrath Mac.pdf'
$ws.Range("A58").Value = 'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf'
$ws.Range("A59").Value = 'Neurological/Sub arachnoid haemorrhage management.pdf'
$ws.Range("A60").Value = 'GI_Liver_and_Transplant/Treatment of constipation.pdf'
$ws.Range("A61").Value = 'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf'
$ws.Range("A63").Value = 'Drugs/anidulafungin.pdf'
$ws.Range("A65").Value = 'GI_Liver_and_Transplant/stress ulcer prophylaxis.pdf'
$ws.Range("A68").Value = 'Drugs/insulin.pdf'
$ws.Range("A69").Value = 'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf'
$ws.Range("A74").Value = 'Post_op_care/Epidural Haematoma.pdf'
$ws.Range("A75").Value = 'Drugs/morphine.pdf'
$ws.Range("A76").Value = 'Breathing(Respiratory)/Equipment/T piece Y piece.pdf'
$ws.Range("A77").Value = 'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf'
$ws.Range("A78").Value = 'Infection_and_sepsis/BAL and MiniBAL standardised procedure.pdf'
$ws.Range("A79").Value = 'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf'
$ws.Range("A80").Value = 'Drugs/atracurium.pdf'
$ws.Range("A87").Value = 'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf'
$ws.Range("A88").Value = 'Drugs/phenylephrine.pdf'
$ws.Range("A89").Value = 'Drugs/amiodarone.pdf'
$ws.Range("A91").Value = 'Drugs/potassium.pdf'
$ws.Range("A92").Value = 'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf'
$ws.Range("A94").Value = 'Drugs/Antibiotic doses in CVVHD.pdf'
$ws.Range("A96").Value = 'Infection_and_sepsis/Antibiotic doses in CVVHD.pdf'
$ws.Range("A98").Value = 'Renal_and_Urology/Antibiotic doses in CVVHD.pdf'
$ws.Range("A100").Value = 'ECLS/RIE ECLS Anti Xa Protocol.pdf'
$ws.Range("A101").Value = 'Drugs/calcium.pdf'
$ws.Range("A102").Value = 'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf'
$ws.Range("A103").Value = 'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf'
$ws.Range("A105").Value = 'Neurological/Management of traumatic brain injury.pdf'
$ws.Range("A106").Value = 'Neurological/Critical Care MRI Procedure_pub.pdf'
$ws.Range("A107").Value = 'Ethics_and_Law/DNACPR policy for Scotland.pdf'
$ws.Range("A108").Value = 'Drugs/vasopressin_sepsis.pdf'
$ws.Range("A109").Value = 'Organ_donation/Organ Retrieval SOP.pdf'
$ws.Range("A110").Value = 'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf'
$ws.Range("A111").Value = 'Drugs/vasopressin organ donation.pdf'
$ws.Range("A112").Value = 'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf'
$ws.Range("A113").Value = 'Ethics_and_Law/Care at the End of Life (FICM).pdf'
$ws.Range("A114").Value = 'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD - Sept 22.pdf'
$ws.Range("A115").Value = 'Drugs/nicardipine.pdf'
$ws.Range("A119").Value = 'Drugs/ketamine for status epilepticus.pdf'
$ws.Range("A120").Value = 'Covid-19/videos/Donning and Doffing Video.pdf'
$ws.Range("A122").Value = 'Drugs/thiopentone.pdf'
$ws.Range("A123").Value = 'Infection_and_sepsis/Infection indications for IVIG.pdf'
$ws.Range("A124").Value = 'Drugs/piperacillin_tazobactam extended_infusion.pdf'
$ws.Range("A125").Value = 'Breathing(Respiratory)/CPAP.pdf'
$ws.Range("A126").Value = 'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf'
$ws.Range("A127").Value = 'Breathing(Respiratory)/Equipment/Bipap V60.pdf'
$ws.Range("A129").Value = 'Covid-19/Covid 19 Death Certification Guideline.pdf'
$ws.Range("A130").Value = 'Transfer/Transfer Outdoors to Garden Guideline.pdf'
$ws.Range("A132").Value = 'Routine_Care/Video Communication.pdf'
$ws.Range("A133").Value = 'Neurological/Treatment of status epilepticus.pdf'
$ws.Range("A135").Value = 'Cardiovascular/Cardiogenic Shock.pdf'
$ws.Range("A136").Value = 'Drugs/isoprenaline.pdf'
$ws.Range("A141").Value = 'Cardiovascular/Management of hypertension within Critical Care.pdf'
$ws.Range("A142").Value = 'Drugs/aminophylline.pdf'
$ws.Range("A143").Value = 'Haematology_CAR-T/CRS.pdf'
$ws.Range("A145").Value = 'Drugs/phenytoin.pdf'
$ws.Range("A147").Value = 'Drugs/pancuronium.pdf'
$ws.Range("A148").Value = 'Drugs/Milrinone.pdf'
$ws.Range("A149").Value = 'Policies_and_admin/General Critical Care SOP_pub.pdf'
$ws.Range("A152").Value = 'Drugs/glyceryl_trinitrate.pdf'
$ws.Range("A155").Value = 'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf'
$ws.Range("A156").Value = 'Drugs/dexmedetomidine.pdf'
$ws.Range("A161").Value = 'Drugs/vancomycin.pdf'
$ws.Range("A162").Value = 'Drugs/neostigmine.pdf'
$ws.Range("A163").Value = 'Drugs/labetalol.pdf'
$ws.Range("A164").Value = 'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf'
$ws.Range("A165").Value = 'Drugs/alteplase for massive PE.pdf'
$ws.Range("A168").Value = 'Routine_Care/ICU Eye Care Guideline.pdf'
$ws.Range("A170").Value = 'Procedures/Arterial Line insertion for ACCPs.pdf'
$ws.Range("A173").Value = 'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf'
$ws.Range("A174").Value = 'Cardiovascular/Cardiac Output Monitoring _pub.pdf'
$ws.Range("A180").Value = 'Procedures/ACCPs acquiring initial CVC competencies.pdf'
$ws.Range("A181").Value = 'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf'
